$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.999.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.753.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.17"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.63"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.750.76"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.99%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.47"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.22%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.90"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.381.73"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.753.73"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.029.81"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.28"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.89"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +19.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.90"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.726"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000153"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.81"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.33"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.38%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.02"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.66"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.898.95"
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.109"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.687.91"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.65%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.87%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.325"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "429.34"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.64"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.31"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.53"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.791.78"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.33%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.44%  "
